# Update as of 2024-04-29
# Adds new "Ingreso" (income) rows and a new "Gastos" (expense) row for the
# date 2024-04-29 (Excel serial date 45411).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Ingreso": append rows 704-711
# ---------------------------------------------------------------------
$wsIngreso = $wb.Worksheets.Item("Ingreso")

$ingresoRows = @(
    @{ Row = 704; Miembro = "Yeyo";      Aporte = 100 },
    @{ Row = 705; Miembro = "Chamo";     Aporte = 200 },
    @{ Row = 706; Miembro = "Carlos";    Aporte = 900 },
    @{ Row = 707; Miembro = "Elio";      Aporte = 100 },
    @{ Row = 708; Miembro = "Invitados"; Aporte = 100 },
    @{ Row = 709; Miembro = "Frandy";    Aporte = 50  },
    @{ Row = 710; Miembro = "Kibelo";    Aporte = 200 },
    @{ Row = 711; Miembro = "Robert";    Aporte = 500 }
)

foreach ($r in $ingresoRows) {
    $row = $r.Row
    $wsIngreso.Cells.Item($row, 1).Value2 = 45411
    $wsIngreso.Cells.Item($row, 2).Value2 = $r.Miembro
    $wsIngreso.Cells.Item($row, 3).Value2 = $r.Aporte
    $wsIngreso.Cells.Item($row, 4).Value2 = "Aporte"
}

# ---------------------------------------------------------------------
# Sheet "Gastos": append row 90
# ---------------------------------------------------------------------
$wsGastos = $wb.Worksheets.Item("Gastos")

$wsGastos.Cells.Item(90, 1).Value2 = 45411
$wsGastos.Cells.Item(90, 2).Value2 = "Arbitro y agua"
$wsGastos.Cells.Item(90, 3).Value2 = 960

# ---------------------------------------------------------------------
# Restore view/selection state as closely as possible
# ---------------------------------------------------------------------
$wsGastos.Activate()
$wsGastos.Range("A90").Select()

$wsIngreso.Activate()
$wsIngreso.Range("C714").Select()
